$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "320.04"); force them to
# stay text (matching the source inline strings) instead of being auto-converted
# to numbers by Excel's smart-entry parsing.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Updated crypto prices / volume deltas
$ws.Range("D2").Value = '30.085.83'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '1.877.26'
$ws.Range("E3").Value = '  -2.23%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '320.04'
$ws.Range("D7").Value = '0.5046'
$ws.Range("D8").Value = '0.3961'
$ws.Range("E8").Value = '  -3.07%  '
$ws.Range("D9").Value = '0.08223'
$ws.Range("E9").Value = '  -3.46%  '
$ws.Range("D10").Value = '42.12'
$ws.Range("E10").Value = '  -1.53%  '
$ws.Range("D11").Value = '1.095'
$ws.Range("E11").Value = '  -2.93%  '
$ws.Range("D12").Value = '23.58'
$ws.Range("E12").Value = '  +5.27%  '
$ws.Range("D13").Value = '1.883.70'
$ws.Range("E13").Value = '  -2.28%  '
$ws.Range("D14").Value = '6.307'
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("E15").Value = '  -2.77%  '
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").Value = '91.93'
$ws.Range("E17").Value = '  -3.99%  '
$ws.Range("D18").Value = '0.00001087'
$ws.Range("E18").Value = '  -2.34%  '
$ws.Range("D19").Value = '0.06487'
$ws.Range("E19").Value = '  -3.00%  '
$ws.Range("D20").Value = '18.14'
$ws.Range("E20").Value = '  -1.29%  '
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").Value = '30.073.49'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").Value = '5.842'
$ws.Range("E23").Value = '  -2.88%  '
$ws.Range("E24").Value = '  -1.54%  '
$ws.Range("D25").Value = '2.167'
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("D26").Value = '2.084.54'
$ws.Range("E26").Value = '  -2.28%  '
$ws.Range("D27").Value = '161.06'
$ws.Range("E27").Value = '  +0.70%  '
$ws.Range("D28").Value = '21.14'
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("E29").Value = '  -8.42%  '
$ws.Range("D30").Value = '127.59'
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("D31").Value = '1.082'
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("D32").Value = '0.1036'
$ws.Range("E32").Value = '  -1.86%  '
$ws.Range("D33").Value = '5.961'
$ws.Range("E33").Value = '  -1.61%  '
$ws.Range("D34").Value = '3.690'
$ws.Range("E34").Value = '  +1.77%  '
$ws.Range("E35").Value = '  -2.32%  '
$ws.Range("D36").Value = '5.287'
$ws.Range("E36").Value = '  +1.71%  '
$ws.Range("D37").Value = '0.06373'
$ws.Range("E37").Value = '  -3.88%  '
$ws.Range("D38").Value = '0.2137'
$ws.Range("E38").Value = '  -3.59%  '
$ws.Range("D39").Value = '1.173'
$ws.Range("E39").Value = '  -5.23%  '
$ws.Range("D40").Value = '8.518'
$ws.Range("E40").Value = '  -4.21%  '
$ws.Range("D41").Value = '0.6300'
$ws.Range("E41").Value = '  -3.87%  '
$ws.Range("D42").Value = '1.216'
$ws.Range("E42").Value = '  -2.67%  '
$ws.Range("D43").Value = '11.30'
$ws.Range("E43").Value = '  -2.93%  '
$ws.Range("D44").Value = '13.25'
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").Value = '0.5912'
$ws.Range("E45").Value = '  -4.10%  '
$ws.Range("D46").Value = '2.094'
$ws.Range("E46").Value = '  +0.39%  '
$ws.Range("D47").Value = '3.632'
$ws.Range("E47").Value = '  -3.67%  '
$ws.Range("D48").Value = '122.38'
$ws.Range("E48").Value = '  -1.88%  '
$ws.Range("D49").Value = '1.209'
$ws.Range("E49").Value = '  -3.37%  '
$ws.Range("D50").Value = '77.56'
$ws.Range("E50").Value = '  -2.84%  '
$ws.Range("D51").Value = '1.114'
$ws.Range("E51").Value = '  -5.26%  '

# Restore default styling on the price column (undo the temporary text format)
$priceRange.Style = "Normal"

